# Update crypto price/volume table to reflect latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price text (e.g. "94.538.97", "0.0425").
# Force text format first so Excel does not coerce these into numbers
# (which would lose trailing zeros / thousands-style dots / precision).
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range("D2").Value = "94.538.97"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "3.532.65"
$ws.Range("E3").Value = "  +6.36%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "238.46"
$ws.Range("E5").Value = "  +4.04%  "
$ws.Range("D6").Value = "633.84"
$ws.Range("E6").Value = "  +2.86%  "
$ws.Range("E7").Value = "  +6.59%  "
$ws.Range("D8").Value = "0.397"
$ws.Range("E8").Value = "  +4.36%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "0.998"
$ws.Range("E10").Value = "  +9.45%  "
$ws.Range("D11").Value = "3.521.71"
$ws.Range("E11").Value = "  +6.16%  "
$ws.Range("D12").Value = "43.64"
$ws.Range("E12").Value = "  +5.53%  "
$ws.Range("E13").Value = "  +5.14%  "
$ws.Range("E14").Value = "  +5.57%  "
$ws.Range("D15").Value = "4.204.77"
$ws.Range("E15").Value = "  +6.43%  "
$ws.Range("D16").Value = "94.333.58"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("E17").Value = "  +4.74%  "
$ws.Range("D18").Value = "8.30"
$ws.Range("E18").Value = "  +5.01%  "
$ws.Range("D19").Value = "3.539.76"
$ws.Range("E19").Value = "  +6.40%  "
$ws.Range("D20").Value = "12.98"
$ws.Range("E20").Value = "  +19.12%  "
$ws.Range("D21").Value = "18.00"
$ws.Range("E21").Value = "  +5.67%  "
$ws.Range("D22").Value = "0.496"
$ws.Range("E22").Value = "  +11.54%  "
$ws.Range("D23").Value = "514.53"
$ws.Range("E23").Value = "  +5.58%  "
$ws.Range("D24").Value = "3.40"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "6.69"
$ws.Range("E25").Value = "  +11.68%  "
$ws.Range("D26").Value = "0.0000193"
$ws.Range("E26").Value = "  +8.13%  "
$ws.Range("D27").Value = "95.84"
$ws.Range("E27").Value = "  +7.49%  "
$ws.Range("D28").Value = "12.22"
$ws.Range("E28").Value = "  +6.55%  "
$ws.Range("D29").Value = "3.03"
$ws.Range("E29").Value = "  +15.53%  "
$ws.Range("D30").Value = "11.55"
$ws.Range("E30").Value = "  +5.57%  "
$ws.Range("D31").Value = "0.143"
$ws.Range("E31").Value = "  +7.08%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "0.182"
$ws.Range("E33").Value = "  +6.71%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").Value = "30.18"
$ws.Range("E35").Value = "  +7.60%  "
$ws.Range("D36").Value = "0.564"
$ws.Range("E36").Value = "  +8.28%  "
$ws.Range("D37").Value = "586.94"
$ws.Range("E37").Value = "  +12.64%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  +7.95%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "7.58"
$ws.Range("E39").Value = "  +4.52%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "0.931"
$ws.Range("E41").Value = "  +7.98%  "
$ws.Range("D42").Value = "0.150"
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "23.76"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0425"
$ws.Range("E44").Value = "  +5.86%  "
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("E46").Value = "  +5.53%  "
$ws.Range("D47").Value = "3.56"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").Value = "2.17"
$ws.Range("E48").Value = "  +4.34%  "
$ws.Range("D49").Value = "53.83"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").Value = "8.13"
$ws.Range("E50").Value = "  +4.67%  "
$ws.Range("E51").Value = "  +3.73%  "
